$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016550561112255
$ws.Range("D2").Value = 1.021986328851226
$ws.Range("E2").Value = 1.018063608546924
$ws.Range("F2").Value = 1.014897203243716
$ws.Range("I2").Value = 1.026413501647625
$ws.Range("J2").Value = 1.021769729978367
$ws.Range("K2").Value = 1.024822041565477
$ws.Range("L2").Value = 1.020910937478898
$ws.Range("M2").Value = 1.017753980080026

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017581810294131
$ws.Range("D3").Value = 1.022721097083471
$ws.Range("E3").Value = 1.018939505950437
$ws.Range("F3").Value = 1.016572484601029
$ws.Range("I3").Value = 1.026583088349896
$ws.Range("J3").Value = 1.022436049147345
$ws.Range("K3").Value = 1.025363576045221
$ws.Range("L3").Value = 1.021592371068638
$ws.Range("M3").Value = 1.019231893363679

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.018248518103569
$ws.Range("D4").Value = 1.023195576963124
$ws.Range("E4").Value = 1.019506175193396
$ws.Range("F4").Value = 1.017655727435359
$ws.Range("I4").Value = 1.026690586933809
$ws.Range("J4").Value = 1.022866089926425
$ws.Range("K4").Value = 1.025712359924329
$ws.Range("L4").Value = 1.022032579390868
$ws.Range("M4").Value = 1.020186985609237

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018528664674083
$ws.Range("D5").Value = 1.023394817455326
$ws.Range("E5").Value = 1.019744380668359
$ws.Range("F5").Value = 1.018110944650347
$ws.Range("I5").Value = 1.026735244049268
$ws.Range("J5").Value = 1.02304661355439
$ws.Range("K5").Value = 1.02585859983071
$ws.Range("L5").Value = 1.022217470038581
$ws.Range("M5").Value = 1.02058822186776

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018575694496977
$ws.Range("D6").Value = 1.023428257273347
$ws.Range("E6").Value = 1.019784375122504
$ws.Range("F6").Value = 1.018187367406278
$ws.Range("I6").Value = 1.026742710781484
$ws.Range("J6").Value = 1.023076908717872
$ws.Range("K6").Value = 1.025883131359133
$ws.Range("L6").Value = 1.022248503867773
$ws.Range("M6").Value = 1.020655574710817

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018252261974015
$ws.Range("D7").Value = 1.023198240131457
$ws.Range("E7").Value = 1.019509358193606
$ws.Range("F7").Value = 1.017661810760183
$ws.Range("I7").Value = 1.026691185747903
$ws.Range("J7").Value = 1.022868503134358
$ws.Range("K7").Value = 1.025714315516745
$ws.Range("L7").Value = 1.022035050587564
$ws.Range("M7").Value = 1.020192348057463

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.01689919735698
$ws.Range("D8").Value = 1.022234846890111
$ws.Range("E8").Value = 1.018359642114088
$ws.Range("F8").Value = 1.01546353782478
$ws.Range("I8").Value = 1.026471277031194
$ws.Range("J8").Value = 1.021995146432122
$ws.Range("K8").Value = 1.025005392316029
$ws.Range("L8").Value = 1.021141381575772
$ws.Range("M8").Value = 1.018253703718146

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.014510427548025
$ws.Range("D9").Value = 1.020529837655088
$ws.Range("E9").Value = 1.016332946642559
$ws.Range("F9").Value = 1.011583606740155
$ws.Range("I9").Value = 1.026066654959197
$ws.Range("J9").Value = 1.020447626157982
$ws.Range("K9").Value = 1.023743720950689
$ws.Range("L9").Value = 1.019561044510838
$ws.Range("M9").Value = 1.01482794673804

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.012914803510839
$ws.Range("D10").Value = 1.019388191088697
$ws.Range("E10").Value = 1.014981278848573
$ws.Range("F10").Value = 1.008992265923486
$ws.Range("I10").Value = 1.025785404709109
$ws.Range("J10").Value = 1.019410141885455
$ws.Range("K10").Value = 1.022894210906287
$ws.Range("L10").Value = 1.018503694968282
$ws.Range("M10").Value = 1.012537219019072

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.012223122311614
$ws.Range("D11").Value = 1.018892662936091
$ws.Range("E11").Value = 1.014395855592136
$ws.Range("F11").Value = 1.007868957908892
$ws.Range("I11").Value = 1.02566089234476
$ws.Range("J11").Value = 1.018959510379517
$ws.Range("K11").Value = 1.022524366814073
$ws.Range("L11").Value = 1.018044941543705
$ws.Range("M11").Value = 1.011543579727029

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.011966084223143
$ws.Range("D12").Value = 1.018708422897209
$ws.Range("E12").Value = 1.014178381103008
$ws.Range("F12").Value = 1.007451515464157
$ws.Range("I12").Value = 1.025614232806781
$ws.Range("J12").Value = 1.018791915140923
$ws.Range("K12").Value = 1.022386689078703
$ws.Range("L12").Value = 1.017874401788302
$ws.Range("M12").Value = 1.01117422833268

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.012021225085965
$ws.Range("D13").Value = 1.018747951109401
$ws.Range("E13").Value = 1.014225031139014
$ws.Range("F13").Value = 1.00754106730044
$ws.Range("I13").Value = 1.025624259995436
$ws.Range("J13").Value = 1.018827874439656
$ws.Range("K13").Value = 1.022416235058812
$ws.Range("L13").Value = 1.017910989419304
$ws.Range("M13").Value = 1.011253467789855

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.012201877858661
$ws.Range("D14").Value = 1.018877437247197
$ws.Range("E14").Value = 1.014377879541495
$ws.Range("F14").Value = 1.007834456036226
$ws.Range("I14").Value = 1.02565704381925
$ws.Range("J14").Value = 1.018945661208704
$ws.Range("K14").Value = 1.022512992476882
$ws.Range("L14").Value = 1.018030847501463
$ws.Range("M14").Value = 1.011513054571563

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.012313168447659
$ws.Range("D15").Value = 1.018957194222939
$ws.Range("E15").Value = 1.014472051515566
$ws.Range("F15").Value = 1.008015196284749
$ws.Range("I15").Value = 1.025677188670053
$ws.Range("J15").Value = 1.019018205577275
$ws.Range("K15").Value = 1.02257256798903
$ws.Range("L15").Value = 1.01810467766512
$ws.Range("M15").Value = 1.011672958593665

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01296069182564
$ws.Range("D16").Value = 1.019421052605098
$ws.Range("E16").Value = 1.015020128393121
$ws.Range("F16").Value = 1.009066789236184
$ws.Range("I16").Value = 1.025793610678253
$ws.Range("J16").Value = 1.019440019342677
$ws.Range("K16").Value = 1.02291871403885
$ws.Range("L16").Value = 1.018534121580468
$ws.Range("M16").Value = 1.012603126243848

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.013366659868046
$ws.Range("D17").Value = 1.01971170047723
$ws.Range("E17").Value = 1.015363883867823
$ws.Range("F17").Value = 1.009726086960964
$ws.Range("I17").Value = 1.025865908406139
$ws.Range("J17").Value = 1.019704237923731
$ws.Range("K17").Value = 1.023135306266721
$ws.Range("L17").Value = 1.018803254990652
$ws.Range("M17").Value = 1.013186124434973

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.013603380432992
$ws.Range("D18").Value = 1.019881115900841
$ws.Range("E18").Value = 1.015564376962499
$ws.Range("F18").Value = 1.010110525139184
$ws.Range("I18").Value = 1.025907815130884
$ws.Range("J18").Value = 1.0198582177374
$ws.Range("K18").Value = 1.023261447858933
$ws.Range("L18").Value = 1.018960147826329
$ws.Range("M18").Value = 1.013526010181404

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.013684083579493
$ws.Range("D19").Value = 1.019938862761267
$ws.Range("E19").Value = 1.015632737635366
$ws.Range("F19").Value = 1.010241588719281
$ws.Range("I19").Value = 1.025922059579877
$ws.Range("J19").Value = 1.019910698111419
$ws.Range("K19").Value = 1.023304426165083
$ws.Range("L19").Value = 1.019013629309441
$ws.Range("M19").Value = 1.013641874319767

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.013323110979638
$ws.Range("D20").Value = 1.019680528564485
$ws.Range("E20").Value = 1.015327003578375
$ws.Range("F20").Value = 1.009655362937209
$ws.Range("I20").Value = 1.025858178781135
$ws.Range("J20").Value = 1.019675903683071
$ws.Range("K20").Value = 1.023112087937274
$ws.Range("L20").Value = 1.018774388646132
$ws.Range("M20").Value = 1.013123591605938

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.012148683361902
$ws.Range("D21").Value = 1.018839311746499
$ws.Range("E21").Value = 1.014332870117191
$ws.Range("F21").Value = 1.007748065848509
$ws.Range("I21").Value = 1.025647401121766
$ws.Range("J21").Value = 1.018910981765252
$ws.Range("K21").Value = 1.022484508147617
$ws.Range("L21").Value = 1.017995556109538
$ws.Range("M21").Value = 1.011436620240929

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.011409596840704
$ws.Range("D22").Value = 1.018309370541171
$ws.Range("E22").Value = 1.013707689919243
$ws.Range("F22").Value = 1.006547734695103
$ws.Range("I22").Value = 1.025512504069873
$ws.Range("J22").Value = 1.01842882572756
$ws.Range("K22").Value = 1.022088181023993
$ws.Range("L22").Value = 1.017505072481073
$ws.Range("M22").Value = 1.010374390907017

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.011801465424231
$ws.Range("D23").Value = 1.018590400624173
$ws.Range("E23").Value = 1.014039122296057
$ws.Range("F23").Value = 1.007184163984253
$ws.Range("I23").Value = 1.025584240490401
$ws.Range("J23").Value = 1.018684541746185
$ws.Range("K23").Value = 1.022298446934099
$ws.Range("L23").Value = 1.017765163360634
$ws.Range("M23").Value = 1.010937649743917

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.013342789075603
$ws.Range("D24").Value = 1.019694614164196
$ws.Range("E24").Value = 1.01534366823599
$ws.Range("F24").Value = 1.009687320446261
$ws.Range("I24").Value = 1.025861672279629
$ws.Range("J24").Value = 1.019688707122139
$ws.Range("K24").Value = 1.023122579897316
$ws.Range("L24").Value = 1.018787432377043
$ws.Range("M24").Value = 1.01315184801574

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.015128523637385
$ws.Range("D25").Value = 1.020971499469631
$ws.Range("E25").Value = 1.016856988172088
$ws.Range("F25").Value = 1.012587457853684
$ws.Range("I25").Value = 1.026173287107494
$ws.Range("J25").Value = 1.020848716315921
$ws.Range("K25").Value = 1.024071371538843
$ws.Range("L25").Value = 1.019970264721034
$ws.Range("M25").Value = 1.015714770205687
